$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.952.12'
$ws.Range('E2').Value = '  -6.70%  '

$ws.Range('D3').Value = '2.178.55'
$ws.Range('E3').Value = '  -7.58%  '

$ws.Range('E4').Value = '  -0.30%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.36'
$ws.Range('E5').Value = '  -0.32%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.617'
$ws.Range('E6').Value = '  -7.69%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '69.35'
$ws.Range('E7').Value = '  -5.46%  '

$ws.Range('E8').Value = '  +0.12%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.533'
$ws.Range('E9').Value = '  -12.60%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '57.46'
$ws.Range('E10').Value = '  -5.29%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '36.06'
$ws.Range('E11').Value = '  +6.57%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0930'
$ws.Range('E12').Value = '  -8.49%  '

$ws.Range('E13').Value = '  -4.27%  '

$ws.Range('E14').Value = '  -10.06%  '

$ws.Range('D15').Value = '2.493.69'
$ws.Range('E15').Value = '  -7.93%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.47'
$ws.Range('E16').Value = '  -10.72%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.826'
$ws.Range('E17').Value = '  -9.24%  '

$ws.Range('D18').Value = '2.186.11'
$ws.Range('E18').Value = '  -7.17%  '

$ws.Range('D19').Value = '40.804.97'
$ws.Range('E19').Value = '  -7.11%  '

$ws.Range('D20').Value = '0.0₃0930'
$ws.Range('E20').Value = '  -9.70%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.88'
$ws.Range('E21').Value = '  -6.26%  '

$ws.Range('E22').Value = '  -8.67%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '229.45'
$ws.Range('E23').Value = '  -9.37%  '

$ws.Range('E24').Value = '  +7.31%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.998'
$ws.Range('E25').Value = '  -0.12%  '

$ws.Range('E26').Value = '  -4.88%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.38'
$ws.Range('E27').Value = '  -4.57%  '

$ws.Range('E28').Value = '  -5.19%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.61'
$ws.Range('E29').Value = '  -8.06%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '166.71'
$ws.Range('E30').Value = '  -5.41%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.00'
$ws.Range('E31').Value = '  -10.18%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.116'
$ws.Range('E32').Value = '  -9.64%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.122'
$ws.Range('E33').Value = '  -8.38%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0694'
$ws.Range('E34').Value = '  -7.00%  '

$ws.Range('E35').Value = '  -5.81%  '

$ws.Range('E36').Value = '  -10.73%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.82'
$ws.Range('E37').Value = '  +0.73%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.30'
$ws.Range('E38').Value = '  +17.72%  '

$ws.Range('E39').Value = '  -7.51%  '

$ws.Range('E40').Value = '  -3.79%  '

$ws.Range('E41').Value = '  -13.13%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '64.61'
$ws.Range('E42').Value = '  -0.95%  '

$ws.Range('B43').Value = 'FTXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.78'
$ws.Range('E43').Value = '  -12.63%  '

$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.60'
$ws.Range('E44').Value = '  -5.05%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.188'
$ws.Range('E45').Value = '  -6.21%  '

$ws.Range('E46').Value = '  -0.05%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0975'
$ws.Range('E47').Value = '  -8.66%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.48'

$ws.Range('B49').Value = 'Celestia'
$ws.Range('C49').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.89'
$ws.Range('E49').Value = '  +3.36%  '

$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.15'
$ws.Range('E50').Value = '  -7.27%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.07'
$ws.Range('E51').Value = '  -7.00%  '
